# RDCC-3540 Upload Template Improvements
# - Rename "Case Worker Data" sheet to "Staff Data"
# - Rename the 8 "Area of WorkN" / "Area of WorkN ID" header columns on that
#   sheet to "ServiceN" / "ServiceN ID"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Case Worker Data")

# Header row column pairs: visible "Area of WorkN" column + hidden "Area of
# WorkN ID" column immediately to its right, for N = 1..8.
$visibleCols = @("M", "O", "Q", "S", "U", "W", "Y", "AA")
$hiddenCols  = @("N", "P", "R", "T", "V", "X", "Z", "AB")

for ($i = 0; $i -lt 8; $i++) {
    $n = $i + 1
    $ws.Range($visibleCols[$i] + "1").Value = "Service" + $n
    $ws.Range($hiddenCols[$i] + "1").Value = "Service" + $n + " ID"
}

# Rename the sheet itself last so the lookup above still resolves by the old name.
$ws.Name = "Staff Data"
